# "Added Control System Plot"
#
# Underlying data edit: update four aircraft-attribute values on the
# "Attributes" sheet (Iyy, Izz, clr, cnb respectively, columns E/F/L/N of
# row 2) and move the sheet's UI selection/scroll position from D7 to L3
# (topLeftCell G1) the way the author last left the worksheet after
# plotting the control-system response.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Attributes")

# --- data changes -----------------------------------------------------
$ws.Range("E2").Value = 35100
$ws.Range("F2").Value = 39600
$ws.Range("L2").Value = 0.03
$ws.Range("N2").Value = 0.1

# --- view/selection changes --------------------------------------------
# Scroll the window so column G is left-most visible (best-effort; mirrors
# the source's <sheetView topLeftCell="G1">) and move the active selection
# to L3 (mirrors <selection activeCell="L3" sqref="L3"/>).
$excel.ActiveWindow.ScrollColumn = 7
$excel.ActiveWindow.ScrollRow = 1
[void]$ws.Range("L3").Select()
